{"js": "// Helper: find the first range matching `text` (exact, case-sensitive) and\n// replace its contents with `replacement`. Throws if not found so mistakes\n// are easy to notice while testing.\nasync function replaceOnce(body, text, replacement) {\n  const results = body.search(text, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + text);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. \"Written: Summer Term 2018\" -> \"Re-Written: Spring Term 2020\"\nawait replaceOnce(\n  body,\n  \"Written: Summer Term 2018\",\n  \"Re-Written: Spring Term 2020\"\n);\n\n// 2. \"Date of Next review: Autumn Term 2019\" -> \"Date of Next review: Spring Term 2021\"\nawait replaceOnce(\n  body,\n  \"Date of Next review: Autumn Term 2019\",\n  \"Date of Next review: Spring Term 2021\"\n);\n\n// 3. \"Proprietor: Ellis Wells / 12th July 2019\" -> \"...17th February 2020\"\nawait replaceOnce(body, \"12th July 2019\", \"17th February 2020\");\n\n// 4. \"the doors are opened and ceases\" -> \"the doors are unlocked and ceases\"\nawait replaceOnce(\n  body,\n  \"the doors are opened and ceases\",\n  \"the doors are unlocked and ceases\"\n);\n\n// 5. \"$20:00 per five minutes\" -> \"\\u00a320:00 per five minutes\" (fix currency symbol)\nawait replaceOnce(\n  body,\n  \"plus $20:00 per five minutes\",\n  \"plus \\u00a320:00 per five minutes\"\n);\n\n// 6. \"related directly\" -> \"referred directly\"\nawait replaceOnce(\n  body,\n  \"All attendance issues should be related directly\",\n  \"All attendance issues should be referred directly\"\n);\n\n// 7. append email address before the final period\nawait replaceOnce(\n  body,\n  \"email our admissions email account.\",\n  \"email our admissions email account currently ewells@wribbenhallschool.co.uk.\"\n);\n\n// 8. add missing closing parenthesis\nawait replaceOnce(\n  body,\n  \"(se \\u201cHolidays\\u201d below. \",\n  \"(se \\u201cHolidays\\u201d below). \"\n);\n\n// 9. \"Care Plan/ is looked after/...\" -> \"Care Plan is looked after/post looked after/...\"\nawait replaceOnce(\n  body,\n  \" Care Plan/ is looked after/a child in need/on the Child Protection Register\",\n  \" Care Plan is looked after/post looked after/a child in need/on the Child Protection Register\"\n);\n\n// 10. \"The pupil is given a mentor, (although\" -> \"The pupil may be given a mentor, (although\"\nawait replaceOnce(\n  body,\n  \"The pupil is given a mentor, (although\",\n  \"The pupil may be given a mentor, (although\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction ReplaceOnce($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. \"Written: Summer Term 2018\" -> \"Re-Written: Spring Term 2020\"\nReplaceOnce \"Written: Summer Term 2018\" \"Re-Written: Spring Term 2020\"\n\n# 2. \"Date of Next review: Autumn Term 2019\" -> \"Date of Next review: Spring Term 2021\"\nReplaceOnce \"Date of Next review: Autumn Term 2019\" \"Date of Next review: Spring Term 2021\"\n\n# 3. \"12th July 2019\" -> \"17th February 2020\"\nReplaceOnce \"12th July 2019\" \"17th February 2020\"\n\n# 4. \"the doors are opened and ceases\" -> \"the doors are unlocked and ceases\"\nReplaceOnce \"the doors are opened and ceases\" \"the doors are unlocked and ceases\"\n\n# 5. \"$20:00 per five minutes\" -> \"\u00a320:00 per five minutes\" (fix currency symbol typo)\nReplaceOnce \"plus `$20:00 per five minutes\" \"plus \u00a320:00 per five minutes\"\n\n# 6. \"related directly\" -> \"referred directly\"\nReplaceOnce \"All attendance issues should be related directly\" \"All attendance issues should be referred directly\"\n\n# 7. append email address before the final period\nReplaceOnce \"email our admissions email account.\" \"email our admissions email account currently ewells@wribbenhallschool.co.uk.\"\n\n# 8. add missing closing parenthesis\nReplaceOnce \"(se \u201cHolidays\u201d below. \" \"(se \u201cHolidays\u201d below). \"\n\n# 9. \"Care Plan/ is looked after/...\" -> \"Care Plan is looked after/post looked after/...\"\nReplaceOnce \" Care Plan/ is looked after/a child in need/on the Child Protection Register\" \" Care Plan is looked after/post looked after/a child in need/on the Child Protection Register\"\n\n# 10. \"The pupil is given a mentor, (although\" -> \"The pupil may be given a mentor, (although\"\nReplaceOnce \"The pupil is given a mentor, (although\" \"The pupil may be given a mentor, (although\"\n"}
